$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8491446345256609
$ws.Range("B3").Value = 0.7288135593220338
$ws.Range("B4").Value = 0.8613013698630136
$ws.Range("B5").Value = 0.3467741935483871
$ws.Range("B6").Value = 0.9691714836223507
$ws.Range("B7").Value = 0.8688472254469468
